$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Orange" item family translations (3 stack sizes share the same en/zh text)
$ws.Range("A4").Value = "ATA_ITEM_ORANGE_3"
$ws.Range("B4").Value = "Orange"
$ws.Range("C4").Value = "橘子"

$ws.Range("A5").Value = "ATA_ITEM_ORANGE_2"
$ws.Range("B5").Value = "Orange"
$ws.Range("C5").Value = "橘子"

$ws.Range("A6").Value = "ATA_ITEM_ORANGE_1"
$ws.Range("B6").Value = "Orange"
$ws.Range("C6").Value = "橘子"

# Eaten orange item + effect text
$ws.Range("A7").Value = "ATA_ITEM_EATEN_ORANGE"
$ws.Range("B7").Value = "Eaten Orange"
$ws.Range("C7").Value = "吃剩的橘子"

$ws.Range("A8").Value = "EFFECT_EATEN_ORANGE"
$ws.Range("B8").Value = "The orange has already been eaten."
$ws.Range("C8").Value = "橘子已经被吃完了"

$ws.Range("A9").Value = "EFFECT_TEMP_ITEM"
$ws.Range("B9").Value = "For the future {0} waves"
$ws.Range("C9").Value = "在接下来的{0}个波次中"

# Widen the en/zh columns to fit the newly added, longer strings.
$ws.Columns.Item(2).ColumnWidth = 32.8
$ws.Columns.Item(3).ColumnWidth = 24.17

$ws.Range("F9").Select() | Out-Null
